# Weekly update: insert a new price report as row 115 for
# "Vega Modelo de Temuco - Rabanito" (La Araucanía), pushing the
# existing rows 115-127 down to 116-128.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 115 (shifts old rows 115..127 -> 116..128)
$ws.Rows.Item(115).Insert()

# Populate the new row 115 with the latest weekly record
$ws.Cells.Item(115, 1).Value  = 10
$ws.Cells.Item(115, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(115, 3).Value  = "La Araucanía"
$ws.Cells.Item(115, 4).Value  = 45194
$ws.Cells.Item(115, 5).Value  = 9
$ws.Cells.Item(115, 6).Value  = 300000001
$ws.Cells.Item(115, 7).Value  = "Rabanito"
$ws.Cells.Item(115, 8).Value  = "Sin especificar"
$ws.Cells.Item(115, 9).Value  = "Primera"
$ws.Cells.Item(115, 10).Value = 80
$ws.Cells.Item(115, 11).Value = 8000
$ws.Cells.Item(115, 12).Value = 8000
$ws.Cells.Item(115, 13).Value = 8000
$ws.Cells.Item(115, 14).Value = "$/docena de paquetes"
$ws.Cells.Item(115, 15).Value = "Provincia de Cautín"
$ws.Cells.Item(115, 16).Value = 667
$ws.Cells.Item(115, 17).Value = 12
$ws.Cells.Item(115, 18).Value = "Hortaliza"
